$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace literal "<br/>" markers with real line breaks in the MSME size-definition cells.
$ws.Range("B22").Value2 = "<200 Manuf., Constr., and Mining, `n<100 Others"
$ws.Range("C23").Value2 = "<NT`$80Millionlion Manuf., Constr., and Mining, `n<NT`$100Millionlion Others"
$ws.Range("B24").Value2 = ">200 Mnf., CnsTurnover, & Minin., `n>100 Others"
$ws.Range("C24").Value2 = ">=NT`$80Millionlion Manuf., Constr., and Mining,`n >=NT`$100Millionlion Others"
